$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the trailing rows (Marlon Brando / Richard Burton / James Cagney) ---
$ws.Rows("6:8").Delete()

# --- Header row (row 1): update existing headers, add new ones ---
$ws.Range("B1").Value = "http://dbpedia.org/ontology/deathPlace"
$ws.Range("C1").Value = "http://dbpedia.org/ontology/parent"

$ws.Range("D1").Value = "http://dbpedia.org/ontology/deathDate"
$ws.Range("E1").Value = "http://dbpedia.org/ontology/birthDate"
$ws.Range("F1").Value = "http://dbpedia.org/ontology/birthPlace"

# Copy the bold/bordered header style from A1 onto the three new header cells
$ws.Range("A1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)

# --- Row 2: Giovanni Francesco Guidi di Bagno ---
$ws.Range("A2").Value = "http://dbpedia.org/resource/Giovanni_Francesco_Guidi_di_Bagno"
$ws.Range("B2").Value = "http://dbpedia.org/resource/Rome"
$ws.Range("C2").Value = "http://dbpedia.org/resource/Colonna_family"
$ws.Range("D2").Value = "http://dbpedia.org/resource/1641"
$ws.Range("E2").Value = "'"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'"
$ws.Range("F2").Style = "Normal"

# --- Row 3: Giovanni Doria ---
$ws.Range("A3").Value = "http://dbpedia.org/resource/Giovanni_Doria"
$ws.Range("B3").Value = "http://dbpedia.org/resource/Palermo"
$ws.Range("C3").Value = "http://dbpedia.org/resource/Giovanni_Andrea_Doria"
$ws.Range("D3").Value = "http://dbpedia.org/resource/1642"
$ws.Range("E3").Value = "'"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'"
$ws.Range("F3").Style = "Normal"

# --- Row 4: Dick Sheppard (priest); B4/C4 stay as they were (already blank) ---
$ws.Range("A4").Value = "http://dbpedia.org/resource/Dick_Sheppard_(priest)"
$ws.Range("D4").Value = "'"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "http://dbpedia.org/resource/1880"
$ws.Range("F4").Value = "http://dbpedia.org/resource/Windsor"

# --- Row 5: Claus Westermann; B5/C5 stay as they were (already blank) ---
$ws.Range("A5").Value = "http://dbpedia.org/resource/Claus_Westermann"
$ws.Range("D5").Value = "'"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "http://dbpedia.org/resource/1909"
$ws.Range("F5").Value = "http://dbpedia.org/resource/Berlin"
